$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Range("S2")
$r.Value = 1
$r.HorizontalAlignment = -4108
Write-Host "done"
